$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.924.22'
$ws.Range('E2').Value = '  -1.91%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.867.10'
$ws.Range('E3').Value = '  -2.54%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.85'
$ws.Range('E5').Value = '  -1.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4992'
$ws.Range('E7').Value = '  -2.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3824'
$ws.Range('E8').Value = '  -3.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08885'
$ws.Range('E9').Value = '  -8.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.118'
$ws.Range('E10').Value = '  -2.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.52'
$ws.Range('E11').Value = '  -1.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.381'
$ws.Range('E12').Value = '  -1.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.66'
$ws.Range('E13').Value = '  -1.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.859.19'
$ws.Range('E14').Value = '  -3.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.233'
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001097'
$ws.Range('E17').Value = '  -3.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '91.06'
$ws.Range('E18').Value = '  -3.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06676'
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.98'
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.105'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.956.26'
$ws.Range('E23').Value = '  -2.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.48'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.285'
$ws.Range('E25').Value = '  -1.44%  '
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.060.94'
$ws.Range('E26').Value = '  -3.61%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.497'
$ws.Range('E27').Value = '  -6.35%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '158.04'
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.64'
$ws.Range('E29').Value = '  -2.59%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.96'
$ws.Range('E30').Value = '  -1.93%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1058'
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.053'
$ws.Range('E32').Value = '  -4.33%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.587'
$ws.Range('E33').Value = '  -1.81%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.599'
$ws.Range('E34').Value = '  -1.13%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.346'
$ws.Range('E35').Value = '  -4.65%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06536'
$ws.Range('E36').Value = '  -2.51%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02393'
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2184'
$ws.Range('E38').Value = '  -1.58%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.278'
$ws.Range('E39').Value = '  +5.89%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.197'
$ws.Range('E40').Value = '  -4.76%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.52'
$ws.Range('E41').Value = '  -1.18%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6366'
$ws.Range('E42').Value = '  -1.05%  '
$ws.Range('B43').Value = 'InternetComputer(DFINITY)'
$ws.Range('C43').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.913'
$ws.Range('E43').Value = '  -2.83%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.16'
$ws.Range('E45').Value = '  -2.80%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5993'
$ws.Range('E46').Value = '  -1.31%  '
$ws.Range('B47').Value = 'WEMIXTOKEN'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.285'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.675'
$ws.Range('E48').Value = '  -2.83%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.224'
$ws.Range('E49').Value = '  +2.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.986'
$ws.Range('E50').Value = '  -3.50%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.72'
$ws.Range('E51').Value = '  -3.26%  '
